# priming.xlsx edit: rename condition column, recode formula results as
# text labels ("credit_card"/"cash") instead of numeric codes (1/2), widen
# column A to fit the new header/labels, and update the sheet selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CC-Helping")
$ws.Activate()

# 1) Header rename: DL_Cond -> CashOrCredit_Ferret
$ws.Range("A1").Value = "CashOrCredit_Ferret"

# 2) Recode the condition formula from numeric (1/2) to text labels.
#    A2 holds a standalone formula; A3:A66 and A67:A95 are each a shared
#    formula group (anchored at A3 and A67 respectively) in the source
#    file, so we re-create that same grouping by setting each range in one
#    shot (keeps row 2 separate from the two shared blocks).
$ws.Range("A2").Formula = '=IF(ISODD(B2),"credit_card", "cash")'
$ws.Range("A3:A66").Formula = '=IF(ISODD(B3),"credit_card", "cash")'
$ws.Range("A67:A95").Formula = '=IF(ISODD(B67),"credit_card", "cash")'

# 3) Widen column A (now holds readable text labels instead of 1/2). This
#    also splits the previously-combined A:B <col> entry in two, leaving
#    column B at its original (default) width.
$ws.Columns("A:A").ColumnWidth = 20

# 4) Scroll the frozen pane back to the top and select column B.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$ws.Columns("B:B").Select()
